# "Final few changes, small additions to analysis of data."
#
# The "group" column (C) used the placeholder codes "A"/"B"; replace them
# with the real condition names used in the dissertation analysis:
#   A -> Placebo
#   B -> Caffeine

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:C64")
$rng.Replace("A", "Placebo", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$rng.Replace("B", "Caffeine", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

# Selection/scroll moved from G6 to G7 (and the view scrolled back to the top).
$ws.Range("G7").Select()
